$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two species observation records that lived in row 2 and row 3 were
# swapped: the "Tretåig hackspett" (Picoides tridactylus) record that used to
# be in row 2 is now in row 3, and the "Trådticka" (Climacocystis borealis)
# record that used to be in row 3 is now in row 2. All the other columns
# (C, P..Z, AA..AE, AG..AY, ...) hold identical values for both rows, so only
# the columns that actually carry per-species data need to move.

$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $tmp = $ws.Range($addr2).Value2
    $ws.Range($addr2).Value2 = $ws.Range($addr3).Value2
    $ws.Range($addr3).Value2 = $tmp
}

# Row 2 (now "Trådticka") has blank placeholder cells in columns J and AF,
# instead of the L/M placeholders (M holding "färska spår") that belonged to
# the "Tretåig hackspett" record now living in row 3. Recreate those blank
# placeholder cells in their new rows (copying a known-blank cell preserves
# an explicit, present-but-empty cell rather than leaving it out entirely)
# and move the "färska spår" activity text to M3.
$ws.Range("I2").Copy($ws.Range("J2"))
$ws.Range("I2").Copy($ws.Range("AF2"))

$ws.Range("I3").Copy($ws.Range("L3"))
$ws.Range("I3").Copy($ws.Range("M3"))
$ws.Range("M3").Value2 = "färska spår"

$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("AF3").ClearContents()
